# Commit: Added Rep4 for Plate 4.
# Adds the "Rep4" (column F) measurement for the Plate 4 rows (271-363)
# of the Microplate sheet, and updates the saved view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Microplate")

$repValues = @{
    271 = 83.06513022273967
    272 = 82.5870797599944
    273 = 56.26339624973034
    274 = 86.82215934586003
    275 = 66.64942405238963
    276 = 63.86237521548974
    277 = 50.02011414452855
    278 = 50.46158798194049
    279 = 58.57781904887243
    280 = 45.33210840049312
    281 = 70.60118331114215
    282 = 52.2504257610636
    283 = 89.13639507058954
    284 = 79.34715611929389
    285 = 45.75060614588172
    286 = 68.02486364714248
    287 = 18.98948822818515
    288 = 13.56969342304095
    289 = 69.131080256622
    290 = 26.09679026936282
    291 = 2.699147826596303
    292 = 64.4670073333677
    293 = 17.32354606549376
    294 = 21.1044646738965
    295 = 60.91868962950837
    296 = 45.1591033843231
    297 = 16.9191844663104
    298 = 18.373507206733
    299 = 8.000124206055535
    300 = 16.45688543746371
    301 = 14.91522118811528
    302 = 5.403032873085098
    303 = 7.749378310874984
    304 = 8.12917628562739
    305 = 23.26864124701497
    306 = 15.08838556028384
    307 = 89.88257452961109
    308 = 32.219401258812
    309 = 11.21348331748481
    310 = 13.62745948965553
    311 = 59.51451993418731
    312 = 42.85580830200518
    313 = 11.87765403184908
    314 = 13.52593584467735
    315 = 12.39404573053123
    316 = 14.69466976978117
    317 = 4.781769860293389
    318 = 14.98872097890112
    319 = 61.42848512181418
    320 = 54.14683809654464
    321 = 6.037470323235544
    322 = 5.319444327658982
    323 = 7.92221732778076
    324 = 50.5976117553884
    325 = 19.80230871036359
    326 = 13.92315756517947
    327 = 17.80934173695986
    328 = 7.376527007392849
    329 = 22.82435073904491
    330 = 43.79331977208233
    331 = 57.71781249338055
    332 = 42.21038310813216
    333 = 13.37028697717174
    334 = 15.98009220058631
    335 = 49.39853616938379
    336 = 3.776752818029277
    337 = 34.33577445670054
    338 = 35.34153935240884
    339 = 15.99175802507365
    340 = 15.06630775677409
    341 = 37.33975239801118
    342 = 41.05138149998693
    343 = 59.90881349501809
    344 = 54.95563564962634
    345 = 16.21330793712889
    346 = 16.87853504735834
    347 = 23.08785658703495
    348 = 28.96326294484966
    349 = 10.35562627611795
    350 = 15.52656798360014
    351 = 37.43790561385918
    352 = 22.17742670570441
    353 = 28.06589067500384
    354 = 46.20600972190894
    355 = 68.98309225802664
    356 = 54.5112501830587
    357 = 89.88653634064124
    358 = 48.58238491276602
    359 = 62.48615412119842
    360 = 21.18416055132712
    361 = 82.17940337186278
    362 = 90.20177902321576
    363 = 84.82448426190808
}

foreach ($row in $repValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $repValues[$row]
}

# Restore the scroll position / active selection recorded in the file
$excel.ActiveWindow.ScrollRow = 265
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G274").Select()
